# Update the GitHub link on the "Demo" slide (Subtitle 2 placeholder):
#   https://github.com/bibinphilip/he-air-pollution/
# becomes
#   https://github.com/bibinphilip/hackerearth_open_innovation/
# split across two runs (".../hackerearth_open_innovation" + "/"), both
# keeping the original rId1 hyperlink + ppaction://hlinkfile action, and
# picking up the same (empty) tooltip already used by the Youtube link
# in the paragraph below.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(11)
$sh = $s.Shapes.Item("Subtitle 2")
$tr = $sh.TextFrame.TextRange

$para1 = $tr.Paragraphs(1)

$prefix    = "GitHub: "
$oldUrl    = "https://github.com/bibinphilip/he-air-pollution/"
$newUrlNoSlash = "https://github.com/bibinphilip/hackerearth_open_innovation"

$linkStart = $para1.Start + $prefix.Length

# Re-text everything up to (but not including) the trailing "/" - this
# rewrites the first run's text and, because it only covers part of the
# original run, splits off a second run for the remaining "/" character
# (which inherits the same rPr / hlinkClick as the first run).
$firstRunLen = $oldUrl.Length - 1
$firstRun = $tr.Characters($linkStart, $firstRunLen)
$firstRun.Text = $newUrlNoSlash

# The range handle above still reports its pre-edit length, so re-fetch
# both resulting runs fresh (by their real post-edit extents) before
# touching their hyperlinks.
$firstRun = $tr.Characters($linkStart, $newUrlNoSlash.Length)
$firstRun.ActionSettings(1).Hyperlink.ScreenTip = ""

$slashStart = $linkStart + $newUrlNoSlash.Length
$secondRun = $tr.Characters($slashStart, 1)
$secondRun.ActionSettings(1).Hyperlink.ScreenTip = ""
